$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.904.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.70%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0884'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.868.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.636.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("E14").Value = '  -0.93%  '
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.923.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.400.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("E36").Value = '  +1.20%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0171'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.776.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("E47").Value = '  -2.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("E49").Value = '  +2.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0505'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.05%  '
